$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 8.656069925401464

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 4.327115817150455

$ws.Range("B4").Value = 0.6545652718822623
$ws.Range("C4").Value = 0.04103571897497393
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 1.378593816457939

$ws.Range("B5").Value = 0.6545652718822623
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 3.223369029078222
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 6.038307959104277

$ws.Range("B6").Value = 3.272327238179451
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 3.223369029078222
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("G6").Value = 8.656069925401464
